$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Replace-ParagraphRangeXml($startPara, $endPara, $xmlBody) {
    $p1 = $d.Paragraphs($startPara)
    $p2 = $d.Paragraphs($endPara)
    $r = $d.Range($p1.Range.Start, $p2.Range.End)
    [void]$r.InsertXML("<w:p $wns>" + $xmlBody)
}

# ---------------------------------------------------------------------
# 1) "Dynamic" / "Filter" paragraphs: move <w:lastRenderedPageBreak/>
#    from the "Filter" run onto the "Dynamic" run.
# ---------------------------------------------------------------------
$xml = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Dynamic</w:t></w:r></w:p>' +
       '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Filter</w:t></w:r></w:p>'
Replace-ParagraphRangeXml 141 142 $xml

# ---------------------------------------------------------------------
# 2) "text-size" / "text-fill" paragraphs: move <w:lastRenderedPageBreak/>
#    from the "text-fill" run onto the "text-size" run.
# ---------------------------------------------------------------------
$xml = '<w:pPr><w:ind w:left="1080" w:firstLine="720"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">  text-size: </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>10;</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>' +
       '<w:p><w:pPr><w:ind w:left="1080" w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">  text-fill: #</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>FFFFFF;</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
Replace-ParagraphRangeXml 93 94 $xml

# ---------------------------------------------------------------------
# 3) "SELECT * FROM ..." / "WHERE ..." paragraphs: move
#    <w:lastRenderedPageBreak/> from the "WHERE" run onto the
#    "SELECT * FROM" run.
# ---------------------------------------------------------------------
$q = [char]39
$xml = '<w:pPr><w:ind w:left="1440" w:firstLine="720"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">SELECT * FROM </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>dcptransportation.citibikestation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/></w:p>' +
       ('<w:p><w:pPr><w:ind w:left="1440" w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">WHERE </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>boro</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>=' + $q + 'M' + $q + '</w:t></w:r></w:p>')
Replace-ParagraphRangeXml 47 48 $xml

# ---------------------------------------------------------------------
# 4) "Data editing" paragraph: strip the yellow highlight, and add a new
#    "Data share" bullet right after it (same list level/style).
# ---------------------------------------------------------------------
$xml = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Data e</w:t></w:r><w:r><w:t>diting</w:t></w:r></w:p>' +
       '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="3"/></w:numPr></w:pPr><w:r><w:t>Data share</w:t></w:r></w:p>'
Replace-ParagraphRangeXml 11 11 $xml
